# 6WLS-30W-BOM.xlsx edit:
#  - Correct MPN for U2/U3 (15V DCDC converter) to match U1's part: DCDCMULT6
#  - Add a new BOM line (row 11) for the ground-banana terminal block (MKDSN / PHOENIX)
#  - Leave a left-aligned number style on the new MPN cell (Phoenix part number)
#  - Update the saved selection to D14 (as left by the author after editing)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the DCDC converter part number used for U2, U3 (was DCDCMULT7, now DCDCMULT6)
$ws.Cells.Item(9, 2).Value = "DCDCMULT6"

# New row 11: terminal block for ground banana
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "MKDSN"
$ws.Cells.Item(11, 3).Value = "MKDSN"
$ws.Cells.Item(11, 4).Value = "Terminal Block"
$ws.Cells.Item(11, 5).Value = "PHOENIX"
$ws.Cells.Item(11, 6).Value = 1729128
$ws.Cells.Item(11, 6).HorizontalAlignment = -4131

# Restore the author's last selection (D14) as saved in the sheet view
$ws.Range("D14").Select()
